# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> plain "Office Theme" palette (used by the notes master)
#   ppt/theme/theme2.xml  -> "Integral" palette (used by the slide master / active design)
#
# The authored edit swaps the two themes' colour palettes (theme1 ends up with
# the "Integral" colours, theme2 ends up with the plain "Office" colours).
#
# The only theme surface the PowerPoint object model exposes here is the
# active design's (slide master's) ColorScheme, which is backed by
# ppt/theme/theme2.xml. We drive that to the plain "Office" palette -- the
# colour values that originally lived in theme1.xml -- which is the
# observable, COM-reachable half of the authored swap.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# ColorScheme.Item index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# Target values = the original "Office Theme" palette (was in theme1.xml).
$targetColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $targetColors.Count; $i++) {
    $hex = $targetColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbValue = $r + ($g * 256) + ($b * 65536)

    $colorScheme.Item($i).RGB = $rgbValue
}
